# Add three new columns (Wins, Losses, Ties) holding the team's season
# record, mirroring the header styling used by the rest of row 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - new columns AC/AD/AE
$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# Match the formatting already used on the rest of the header row (A1 has
# the bold/centered/bordered header style) by copying its format onto the
# freshly written header cells.
$ws.Range("A1").Copy()
$ws.Range("AC1:AE1").PasteSpecial(-4122)

# Data rows (2 through 41) - season record is the same for every player
# row on this sheet: 83 wins, 79 losses, 0 ties.
for ($r = 2; $r -le 41; $r++) {
    $ws.Cells.Item($r, 29).Value = 83
    $ws.Cells.Item($r, 30).Value = 79
    $ws.Cells.Item($r, 31).Value = 0
}
